# Regenerate the handoff report: the e2e test doc's generated GUID/hash
# changed (new handoff), so every reference to the old identifiers needs
# to be swapped for the new ones, and the handoff timestamps bumped.

$wb = $excel.ActiveWorkbook

$oldGuid = "f81075dc-494b-4358-a6aa-77285574366f"
$newGuid = "395cfa96-5738-4dc2-b9a0-bd6a4c420f49"

$oldHash = "2853dfba4656b750e10f28d8b08b5328141bda18"
$newHash = "558ebfd6393f39351e257d678d2c88b4ebbd1304"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a80921db894c6415fd03337ec4e1ef988148910/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B2").Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")

$wsOverview.Range("G2").Value = "2016-10-26 07:54:38"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Range("A2").Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-10-26 07:54:26"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Range("A2").Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-10-26 07:54:38"
